$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "MODEL_CONDITION" header to "MODELCONDITION" (drop the underscore)
$ws.Range("E1").Value = "MODELCONDITION"

# The old column A (a redundant taxon-id helper column) is no longer needed.
# Deleting it shifts the remaining data (old columns B:F) left into A:E.
$ws.Range("A:A").Delete()
